$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10000.5
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = $null
$ws.Range("H33").Value = 496.7
$ws.Range("I33").Value = 607.9231
$ws.Range("K33").Value = 607.9231
$ws.Range("M33").Value = -378.9231
$ws.Range("H40").Value = 3815.8635
$ws.Range("I40").Value = 4243.3335
$ws.Range("J40").Value = 2899.8572
$ws.Range("K40").Value = 4243.3335
$ws.Range("L40").Value = 2899.8572
$ws.Range("M40").Value = -4068.3335
$ws.Range("N40").Value = -3249.8572
$ws.Range("H86").Value = 7374.0557
$ws.Range("I86").Value = 6188.125
$ws.Range("J86").Value = 8322.799999999999
$ws.Range("K86").Value = 6188.125
$ws.Range("L86").Value = 8322.799999999999
$ws.Range("M86").Value = -5065.125
$ws.Range("N86").Value = -10568.8
$ws.Range("H89").Value = 7374.0557
$ws.Range("I89").Value = 6188.125
$ws.Range("J89").Value = 8322.799999999999
$ws.Range("K89").Value = 30940.625
$ws.Range("L89").Value = 41614
$ws.Range("M89").Value = -25324.625
$ws.Range("N89").Value = -52846
$ws.Range("H106").Value = 111112530
$ws.Range("J106").Value = 3499
$ws.Range("L106").Value = 3499
$ws.Range("N106").Value = -4761
$ws.Range("H107").Value = 19607986
$ws.Range("I107").Value = 33333406
$ws.Range("K107").Value = 33333406
$ws.Range("M107").Value = -33331486
$ws.Range("H115").Value = 778.41174
$ws.Range("I115").Value = 764.5625
$ws.Range("J115").Value = 1000
$ws.Range("K115").Value = 2293.6875
$ws.Range("L115").Value = 3000
$ws.Range("M115").Value = -726.6875
$ws.Range("N115").Value = -6134
$ws.Range("H116").Value = 6369.1904
$ws.Range("J116").Value = 7846.273
$ws.Range("L116").Value = 7846.273
$ws.Range("N116").Value = -14730.273
$ws.Range("H131").Value = 5381.893
$ws.Range("I131").Value = 1815.2142
$ws.Range("K131").Value = 5445.642599999999
$ws.Range("M131").Value = -405.6425999999992
$ws.Range("H137").Value = 40868.22
$ws.Range("I137").Value = 52741.715
$ws.Range("K137").Value = 158225.145
$ws.Range("M137").Value = -155675.145
$ws.Range("H138").Value = 2988.194
$ws.Range("I138").Value = 1221.7059
$ws.Range("J138").Value = 3588.8
$ws.Range("K138").Value = 3665.1177
$ws.Range("L138").Value = 10766.4
$ws.Range("M138").Value = 1474.8823
$ws.Range("N138").Value = -21046.4
$ws.Range("H140").Value = 101394.75
$ws.Range("J140").Value = 101394.75
$ws.Range("L140").Value = 101394.75
$ws.Range("N140").Value = -111754.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10885.197
$ws.Range("I32").Value = 6785.156
$ws.Range("J32").Value = 19671
$ws.Range("K32").Value = 6785.156
$ws.Range("L32").Value = 19671
$ws.Range("M32").Value = -6498.156
$ws.Range("N32").Value = -20245
$ws.Range("H102").Value = 8338914.5
$ws.Range("I102").Value = 10421143
$ws.Range("K102").Value = 10421143
$ws.Range("M102").Value = -10419521
$ws.Range("H132").Value = 3175.6785
$ws.Range("I132").Value = 2206.4
$ws.Range("J132").Value = 5598.875
$ws.Range("K132").Value = 6619.200000000001
$ws.Range("L132").Value = 16796.625
$ws.Range("M132").Value = -4089.200000000001
$ws.Range("N132").Value = -21856.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6805542
$ws.Range("I99").Value = 10206174
$ws.Range("J99").Value = 4279.143
$ws.Range("K99").Value = 10206174
$ws.Range("L99").Value = 4279.143
$ws.Range("M99").Value = -10204676
$ws.Range("N99").Value = -7275.143
$ws.Range("H107").Value = 1787636.6
$ws.Range("I107").Value = 2305921.5
$ws.Range("J107").Value = 2432.5557
$ws.Range("K107").Value = 2305921.5
$ws.Range("L107").Value = 2432.5557
$ws.Range("M107").Value = -2304001.5
$ws.Range("N107").Value = -6272.5557
$ws.Range("H109").Value = 89995
$ws.Range("J109").Value = 89995
$ws.Range("L109").Value = 89995
$ws.Range("N109").Value = -92769
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null
$ws.Range("H122").Value = 50001
$ws.Range("J122").Value = 50001
$ws.Range("L122").Value = 50001
$ws.Range("N122").Value = -59801
$ws.Range("H134").Value = 4194.4062
$ws.Range("I134").Value = 2117.5386
$ws.Range("K134").Value = 6352.6158
$ws.Range("M134").Value = -3817.6158
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28966.6
$ws.Range("I31").Value = 1666.2632
$ws.Range("J31").Value = 61385.75
$ws.Range("K31").Value = 1666.2632
$ws.Range("L31").Value = 61385.75
$ws.Range("M31").Value = -1371.2632
$ws.Range("N31").Value = -61975.75
$ws.Range("H34").Value = 28966.6
$ws.Range("I34").Value = 1666.2632
$ws.Range("J34").Value = 61385.75
$ws.Range("K34").Value = 1666.2632
$ws.Range("L34").Value = 61385.75
$ws.Range("M34").Value = -1464.2632
$ws.Range("N34").Value = -61789.75
$ws.Range("H58").Value = 6234.269
$ws.Range("I58").Value = 9179.923000000001
$ws.Range("J58").Value = 3288.6155
$ws.Range("K58").Value = 9179.923000000001
$ws.Range("L58").Value = 3288.6155
$ws.Range("M58").Value = -8976.923000000001
$ws.Range("N58").Value = -3694.6155
$ws.Range("H86").Value = 9287.733
$ws.Range("I86").Value = 5718.7144
$ws.Range("J86").Value = 12410.625
$ws.Range("K86").Value = 5718.7144
$ws.Range("L86").Value = 12410.625
$ws.Range("M86").Value = -4595.7144
$ws.Range("N86").Value = -14656.625
$ws.Range("H89").Value = 9287.733
$ws.Range("I89").Value = 5718.7144
$ws.Range("J89").Value = 12410.625
$ws.Range("K89").Value = 28593.572
$ws.Range("L89").Value = 62053.125
$ws.Range("M89").Value = -22977.572
$ws.Range("N89").Value = -73285.125
$ws.Range("H122").Value = 1912.25
$ws.Range("I122").Value = 1795.6786
$ws.Range("J122").Value = 2320.25
$ws.Range("K122").Value = 5387.0358
$ws.Range("L122").Value = 6960.75
$ws.Range("M122").Value = -2937.0358
$ws.Range("N122").Value = -11860.75
$ws.Range("H132").Value = 65321.168
$ws.Range("I132").Value = 40336.42
$ws.Range("J132").Value = 227722
$ws.Range("K132").Value = 121009.26
$ws.Range("L132").Value = 683166
$ws.Range("M132").Value = -118479.26
$ws.Range("N132").Value = -688226
$ws.Range("H136").Value = 6234.269
$ws.Range("I136").Value = 9179.923000000001
$ws.Range("J136").Value = 3288.6155
$ws.Range("K136").Value = 27539.769
$ws.Range("L136").Value = 9865.8465
$ws.Range("M136").Value = -24989.769
$ws.Range("N136").Value = -14965.8465
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 927.5
$ws.Range("J122").Value = 750.6429000000001
$ws.Range("L122").Value = 6755.7861
$ws.Range("N122").Value = -11655.7861
$ws.Range("H132").Value = 2376.5417
$ws.Range("I132").Value = 1361.7142
$ws.Range("J132").Value = 2794.4119
$ws.Range("K132").Value = 12255.4278
$ws.Range("L132").Value = 25149.7071
$ws.Range("M132").Value = -9725.427799999999
$ws.Range("N132").Value = -30209.7071
$ws.Range("H139").Value = 2065
$ws.Range("I139").Value = 1593.75
$ws.Range("K139").Value = 4781.25
$ws.Range("M139").Value = 358.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10530707
$ws.Range("I70").Value = 20003794
$ws.Range("J70").Value = 5056.4443
$ws.Range("K70").Value = 20003794
$ws.Range("L70").Value = 5056.4443
$ws.Range("M70").Value = -20003524
$ws.Range("N70").Value = -5596.4443
$ws.Range("H73").Value = 10530707
$ws.Range("I73").Value = 20003794
$ws.Range("J73").Value = 5056.4443
$ws.Range("K73").Value = 20003794
$ws.Range("L73").Value = 5056.4443
$ws.Range("M73").Value = -20002858
$ws.Range("N73").Value = -6928.4443
$ws.Range("H126").Value = 3581179.8
$ws.Range("I126").Value = 5052862.5
$ws.Range("J126").Value = 3090619
$ws.Range("K126").Value = 15158587.5
$ws.Range("L126").Value = 9271857
$ws.Range("M126").Value = -15156117.5
$ws.Range("N126").Value = -9276797
$ws.Range("H132").Value = 4707.727
$ws.Range("I132").Value = 4644
$ws.Range("K132").Value = 13932
$ws.Range("M132").Value = -11402
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1094.125
$ws.Range("I68").Value = 958.8333
$ws.Range("K68").Value = 958.8333
$ws.Range("M68").Value = -209.8333
$ws.Range("H71").Value = 1094.125
$ws.Range("I71").Value = 958.8333
$ws.Range("K71").Value = 4794.1665
$ws.Range("M71").Value = -1050.1665
$ws.Range("H132").Value = 12326.533
$ws.Range("I132").Value = 12530.615
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 37591.845
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = -35061.845
$ws.Range("N132").Value = -38060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21064326
$ws.Range("I132").Value = 23260112
$ws.Range("K132").Value = 69780336
$ws.Range("M132").Value = -69777806
$ws.Range("H137").Value = 80519.82000000001
$ws.Range("J137").Value = 80519.82000000001
$ws.Range("L137").Value = 80519.82000000001
$ws.Range("N137").Value = -90719.82000000001
$ws.Range("H139").Value = 156137
$ws.Range("J139").Value = 156137
$ws.Range("L139").Value = 156137
$ws.Range("N139").Value = -166417
